$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a new worksheet "2022-Q3" right after "总计" (i.e. before the
#    current 2nd sheet "2022-Q2"), matching the workbook.xml diff which
#    shifts every existing quarter sheet down by one position and inserts
#    the new quarter as sheetId 2 / position 2.
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item(1)
$q2Sheet = $wb.Worksheets.Item(2)
$newSheet = $wb.Worksheets.Add($q2Sheet)
$newSheet.Name = "2022-Q3"

# ---------------------------------------------------------------------------
# 2. Populate the new "2022-Q3" sheet with the same layout used by the other
#    quarterly fund-holding sheets: header row + one row per fund.
# ---------------------------------------------------------------------------
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")

# Copy the header style (bold + border, style index 2 in the original file)
# from the existing "总计" sheet header cells so the new sheet matches the
# look of the sibling quarter sheets exactly.
$totalSheet.Range("B1:D1").Copy() | Out-Null
$newSheet.Range("B1:H1").PasteSpecial(-4122) | Out-Null
$newSheet.Application.CutCopyMode = $false

for ($col = 2; $col -le 8; $col++) {
    $newSheet.Cells.Item(1, $col).Value = $headers[$col - 2]
}

$fundRows = @(
    @("006551", "中庚价值领航混合", "118.19", "91.86", "5.32", "6.2877", 6),
    @("007130", "中庚小盘价值股票", "75.87", "93.06", "5.59", "4.2411", 2),
    @("011174", "中庚价值品质一年持有期混合", "66.33", "92.24", "4.52", "2.9981", 6),
    @("090018", "大成新锐产业混合", "95.85", "83.49", "2.53", "2.4250", 10),
    @("001300", "大成睿景灵活配置混合A", "33.29", "90.96", "2.54", "0.8456", 10),
    @("013435", "大成景气精选六个月持有混合A", "32.65", "88.02", "2.51", "0.8195", 10),
    @("001301", "大成睿景灵活配置混合C", "23.86", "90.96", "2.54", "0.6060", 10),
    @("013436", "大成景气精选六个月持有混合C", "5.65", "88.02", "2.51", "0.1418", 10),
    @("160620", "鹏华中证A股资源产业指数（LOF）A", "1.75", "94.44", "2.15", "0.0376", 10),
    @("260117", "景顺长城支柱产业混合", "0.68", "89.87", "4.50", "0.0306", 8),
    @("006441", "中信建投中证500指数增强C", "2.25", "93.60", "1.15", "0.0259", 4),
    @("006440", "中信建投中证500指数增强A", "2.14", "93.60", "1.15", "0.0246", 4),
    @("013878", "圆信永丰中证500指数增强A", "0.95", "92.59", "1.36", "0.0129", 8),
    @("159990", "银华巨潮小盘价值ETF", "0.80", "96.51", "1.27", "0.0102", 7),
    @("012808", "鹏华中证A股资源产业指数（LOF）C", "0.24", "94.44", "2.15", "0.0052", 10),
    @("002952", "建信多因子量化股票", "0.09", "91.26", "3.73", "0.0034", 4),
    @("013879", "圆信永丰中证500指数增强C", "0.09", "92.59", "1.36", "0.0012", 8)
)

$newSheet.Range("B2:G18").NumberFormat = "@"

$rowIndex = 2
foreach ($fund in $fundRows) {
    $newSheet.Cells.Item($rowIndex, 1).Value = $rowIndex - 2
    $newSheet.Cells.Item($rowIndex, 2).Value = $fund[0]
    $newSheet.Cells.Item($rowIndex, 3).Value = $fund[1]
    $newSheet.Cells.Item($rowIndex, 4).Value = $fund[2]
    $newSheet.Cells.Item($rowIndex, 5).Value = $fund[3]
    $newSheet.Cells.Item($rowIndex, 6).Value = $fund[4]
    $newSheet.Cells.Item($rowIndex, 7).Value = $fund[5]
    $newSheet.Cells.Item($rowIndex, 8).Value = $fund[6]
    $rowIndex = $rowIndex + 1
}

# Match the A-column style (bold, centered index style) used on every other
# quarterly sheet.
$totalSheet.Range("A2").Copy() | Out-Null
$newSheet.Range("A2:A18").PasteSpecial(-4122) | Out-Null
$newSheet.Application.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3. Update the "总计" (summary) sheet: insert a new row for 2022-Q3 right
#    after the header row, pushing the existing quarters down, and append
#    the newly revealed 2020-Q4 row at the bottom.
# ---------------------------------------------------------------------------
$totalSheet.Rows.Item(2).Insert()
$totalSheet.Range("A2:D2").ClearFormats()

$totalSheet.Range("B2").NumberFormat = "@"
$totalSheet.Range("B2").Value = "2022-Q3"
$totalSheet.Range("C2").Value = 17
$totalSheet.Range("D2").Value = 18.52

$totalSheet.Range("B9").NumberFormat = "@"
$totalSheet.Range("B9").Value = "2020-Q4"
$totalSheet.Range("A9").Value = 7
$totalSheet.Range("C9").Value = 4
$totalSheet.Range("D9").Value = 2.73

# Re-apply the index-column style (bold, centered - style index 2 in the
# original workbook) to A2 and A9 so every row in column A looks the same,
# then rewrite the whole column A as a fresh 0-based row index (the
# left-most "序号" column always mirrors the row position in this sheet).
$totalSheet.Range("A3").Copy() | Out-Null
$totalSheet.Range("A2").PasteSpecial(-4122) | Out-Null
$totalSheet.Range("A9").PasteSpecial(-4122) | Out-Null
$totalSheet.Application.CutCopyMode = $false

for ($row = 2; $row -le 9; $row++) {
    $totalSheet.Cells.Item($row, 1).Value = $row - 2
}

# Keep the original active sheet ("总计") selected, exactly as it was
# before the edit.
$totalSheet.Activate()
